$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0, 0.4965871691703796, 0.006825661659240723),
    @(1, 0.7680098652839661, 0.1392850526393727),
    @(2, 0.782154631614685, 0.04232812892036745),
    @(3, 0.7060444116592407, 0.119835785011813),
    @(4, 0.8953124761581421, 0.02778017493619667),
    @(5, 0.9010690808296203, 0.04489960378837834),
    @(6, 0.5525082230567933, 0.05585504138590723),
    @(7, 0.9309621691703797, 0.05136255388408779),
    @(8, 0.9233141303062439, 0.03320143014524354),
    @(9, 0.5225329041481018, 0.01140590219950456),
    @(10, 0.4038926524777252, 0.0541943277279549),
    @(11, 0.7513985658799107, 0.1653469986355771),
    @(12, 0.7808253502520608, 0.04249130391582516),
    @(13, 0.6935019177573688, 0.1295943372390366),
    @(14, 0.8948479310292814, 0.02811452200975032),
    @(15, 0.9002750536130222, 0.04519054474493579),
    @(16, 0.4794554695150628, 0.1035783530436264),
    @(17, 0.9296867932145789, 0.052589137913941),
    @(18, 0.9232541603909981, 0.03320371636364689),
    @(19, 0.4633326887508805, 0.0242152832460079),
    @(20, 0.5004718474073753, 0.003803346132688842),
    @(21, 0.8498151699922092, 0.1500610527910309),
    @(22, 0.8794450839140409, 0.02949058524935592),
    @(23, 0.795468317174515, 0.1648746561059093),
    @(24, 0.9665202048346607, 0.01973715835830648),
    @(25, 0.9703931137465375, 0.0252654316639312),
    @(26, 0.6432569847212604, 0.09029948913893696),
    @(27, 0.9863806050900277, 0.01434377178657925),
    @(28, 0.9822067472082756, 0.01735241245652931),
    @(29, 0.497912090438885, 0.04113016954811693)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    if ($row -gt 21) {
        $ws.Range("A21").Copy()
        $ws.Cells.Item($row, 1).PasteSpecial(-4122)
    }
    $row++
}
